# Generate Report for Handback
#
# The handback run for the 9e45521e-a2d2-475d-85ea-8856077458eb source file has
# now completed for de-de as well (it previously only showed the zh-cn
# handback), so the recorded handoff/handback timestamps for that row move
# forward on both the zh-cn and de-de detail sheets, and the Overview sheet's
# summary cell for that file/locale reflects the new handback time.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets("Overview")
$wsOverview.Range("G3").Value = "2016-08-13 17:05:43"

$wsZhCn = $wb.Worksheets("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-13 17:05:35"
$wsZhCn.Range("K3").Value = "2016-08-13 17:06:07"

$wsDeDe = $wb.Worksheets("de-de")
$wsDeDe.Range("H3").Value = "2016-08-13 17:05:43"
$wsDeDe.Range("K3").Value = "2016-08-13 17:06:17"
